$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "next available ID" range for vaccine adjuvant terms.
# Two new vaccine adjuvant terms were added (VO_0005507, VO_0005508),
# so the starting ID for the reserved range moves from VO_0005507 to VO_0005509.
$ws.Range("A2").Value = "VO_0005509 - VO_0005560"

# Reflect the cell selection left active in the sheet when last saved.
[void]$ws.Range("C5").Select()
